# agregando img_not_available y editando agregar en mostrar prod
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("templates")

# H3: new comment "L: js para ingredientes"
$ws.Range("H3").Value = "L: js para ingredientes"

# Row 4 grows taller to fit the new multi-line comment in H4
$ws.Rows.Item(4).RowHeight = 50.25

# H4: new multi-line comment about no_image_available sizing + filters
$ws.Range("H4").Value = "L: mismo tamaño para no_image_available`nG: filtros (activo, inactivo, todos, tipo)"

# Column H (comentarios) widened and no longer best-fit
$ws.Columns.Item(8).ColumnWidth = 22

# Move the active selection to H5
$ws.Activate() | Out-Null
$ws.Range("H5").Select() | Out-Null
